# JDBC performance data collection
# Adds a new "benchmark" row (19) ahead of the second results table, retitles
# the two lower tables from the ORM-NF labels to the new JDBC-NF labels, and
# replaces the underlying measurement data for both of those tables with the
# freshly collected JDBC numbers. Downstream SUM / weighted-average formulas
# in columns H:M recompute automatically from the edited inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 19: a lone "benchmark" marker row sitting above the JDBC 3NF
#     table, mirroring the one that already precedes the first table (row 1).
$ws.Range("A19").Value = "benchmark"
$ws.Range("E19").Value = 21

# --- Re-title the two summary tables.
$ws.Range("A20").Value = "JDBC 3NF"
$ws.Range("A28").Value = "JDBC 0NF"

# --- JDBC 3NF measurements (rows 21-25); B/F mostly unchanged.
$ws.Range("C21").Value = 702016
$ws.Range("D21").Value = 59.06
$ws.Range("E21").Value = 22.5
$ws.Range("G21").Value = 1392.35

$ws.Range("C22").Value = 959162
$ws.Range("D22").Value = 56.86
$ws.Range("E22").Value = 23.68
$ws.Range("G22").Value = 1417.82

$ws.Range("C23").Value = 1061299
$ws.Range("D23").Value = 52.71
$ws.Range("E23").Value = 24.54
$ws.Range("F23").Value = 0.02
$ws.Range("G23").Value = 1336.32

$ws.Range("C24").Value = 1094690
$ws.Range("D24").Value = 49.31
$ws.Range("E24").Value = 25.13
$ws.Range("F24").Value = 0.13
$ws.Range("G24").Value = 1292.45

$ws.Range("C25").Value = 1120446
$ws.Range("D25").Value = 48.91
$ws.Range("E25").Value = 25.73
$ws.Range("F25").Value = 0.63
$ws.Range("G25").Value = 1237.63

# --- JDBC 0NF measurements (rows 29-32); row 33 keeps its original inputs,
#     only its rolling formulas move because the totals above it changed.
$ws.Range("C29").Value = 312518
$ws.Range("D29").Value = 53.22
$ws.Range("E29").Value = 22.49
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 1979.29

$ws.Range("C30").Value = 471503
$ws.Range("D30").Value = 46.59
$ws.Range("E30").Value = 22.63
$ws.Range("F30").Value = 0.32
$ws.Range("G30").Value = 2069.19

$ws.Range("C31").Value = 495433
$ws.Range("D31").Value = 43.41
$ws.Range("E31").Value = 22.21
$ws.Range("F31").Value = 0.02
$ws.Range("G31").Value = 2091.42

$ws.Range("C32").Value = 488661
$ws.Range("D32").Value = 42.98
$ws.Range("E32").Value = 22.23
$ws.Range("F32").Value = 0.01
$ws.Range("G32").Value = 2073.01

# --- Move the cursor/selection to G25, scrolled so row 8 is visible at top.
$ws.Range("G25").Select()
$w = $excel.ActiveWindow
$w.ScrollRow = 8
$w.ScrollColumn = 1
